$d = $word.ActiveDocument

# The original paragraph 1 ends with a hidden "_GoBack" bookmark. Remove it now;
# we'll re-create it at the very end of the document once all new content is in place
# (mirrors where Word leaves it after the last edit).
$gb = $d.Bookmarks.Item("_GoBack")
$gb.Delete()

# --- Paragraph 1: extend "Documents" ---
$p1 = $d.Paragraphs(1).Range
$p1.Collapse(0)
$p1.InsertAfter(" needs to ")
$p1.Collapse(0)
$p1.InsertAfter("represented as vectors.")

# --- Paragraph 2 (new) ---
$p1.Collapse(0)
$p1.InsertParagraphAfter()
$p2 = $d.Paragraphs(2).Range
$p2.InsertAfter("Ideally ")
$p2.Collapse(0)
$p2.InsertAfter("they should be normalised in respect to the whole collection, for better analysis and easier clustering. ")

# --- Paragraph 3 (new) ---
$p2.Collapse(0)
$p2.InsertParagraphAfter()
$p3 = $d.Paragraphs(3).Range
$p3.InsertAfter("At the moment there is not possibility to just add a new document to the analysis process while analysing and/or clustering. The process needs to be rerun from step 0 performing the whole pipeline. Addi")
$p3.Collapse(0)
$p3.InsertAfter("ng this functionality would involve ")
$p3.Collapse(0)
$p3.InsertAfter("caching the")
$p3.Collapse(0)
$p3.InsertAfter(" process at the current stage and then adding the new document(s) to perform the analysis in regard to this document (very difficult to achieve with analysis involving algorithms such IDF where the vectors are weighted strictly in respect to the others). Perhaps it is achievable with a different ")
$p3.Collapse(0)
$p3.InsertAfter("method of")
$p3.Collapse(0)
$p3.InsertAfter(" clustering, especially flat. The weighted document could be clustered simply on respect already saved values coming from the already existing clusters. When it comes to a similarity hierarchy, it may involve the computation of an entirely new similarity matrix. ")
$p3.Collapse(0)
$p3.InsertAfter(" ")

# --- Paragraph 4 (new, empty) ---
$p3.Collapse(0)
$p3.InsertParagraphAfter()

# --- Paragraph 5 (new) ---
$p4 = $d.Paragraphs(4).Range
$p4.Collapse(0)
$p4.InsertParagraphAfter()
$p5 = $d.Paragraphs(5).Range
$p5.InsertAfter("Only a few type")
$p5.Collapse(0)
$p5.InsertAfter("s")
$p5.Collapse(0)
$p5.InsertAfter(" of document")
$p5.Collapse(0)
$p5.InsertAfter("s")
$p5.Collapse(0)
$p5.InsertAfter(" are available at the moment")
$p5.Collapse(0)
$p5.InsertAfter(": Microsoft Word format (97-200")
$p5.Collapse(0)
$p5.InsertAfter("3 and ")
$p5.Collapse(0)
$p5.InsertAfter("docx")
$p5.Collapse(0)
$p5.InsertAfter("), standard txt and other related. There is some experiment ")

# Re-create the "_GoBack" bookmark at the very end of the document content,
# right after the text we just inserted (mirrors the cursor's last-edit spot).
# A bookmark collapsed exactly at the story's last position gets mis-placed by
# this engine, so: append a one-char placeholder, anchor the bookmark right
# before it (a safe, interior position), then remove the placeholder again.
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertAfter("X")

$xPos = $d.Content.End - 2
$bmRange = $d.Range($xPos, $xPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($xPos, $xPos + 1)
$placeholder.Delete()
